$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Unidades" (row 14) loses its (visually inert) style and reverts to the
# sheet's default formatting.
$ws.Range("A14").ClearFormats()

# The old blank, underlined placeholder line (row 15) becomes a real
# labeled field ("Dirección del Cliente"). Clear its inherited underline
# formatting first so the new label renders with plain/default formatting.
$a15 = $ws.Range("A15")
$a15.ClearFormats()
$a15.Value = "Dirección del Cliente"

# A brand-new labeled field is appended: "Inicio Calibración" (row 16).
$a16 = $ws.Range("A16")
$a16.Value = "Inicio Calibración"

# A new blank input line is appended at the end (row 17), mirroring the
# placeholder row that used to sit at row 15.
$a17 = $ws.Cells.Item(17, 1)
$a17.Font.Underline = $true
$a17.Font.Underline = $false

# Move the selection/active cell to the new last row, matching the saved
# view state of the edited workbook.
$ws.Range("A17").Select()
